$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move forward by one day
# (45204 -> 45205, i.e. 2023-10-05 -> 2023-10-06) for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value2 = 45205
    }
}
